# Bento local search - update StatQuery column (C) on the "startup" sheet
# with the working Neo4j "selection script" (counts Programs/Arms/Cases/
# Samples/Assays/Files), update the active selection, and adjust row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new Cypher query that replaces the old "study_subjects" count query.
# Built with backtick-n so each segment lands on its own line in the cell,
# and ends with a trailing newline just like the authored content.
$newQuery = "MATCH (ss:study_subject)`n" +
            "MATCH (ss)<-[:sf_of_study_subject]-(sf)`n" +
            "MATCH (ss)<-[:diagnosis_of_study_subject]-(d)`n" +
            "MATCH (d)<-[:tp_of_diagnosis]-(tp)`n" +
            "MATCH (ss)-[:study_subject_of_study]->(s)`n" +
            "WHERE ss.study_subject_id = 'BENTO-CASE-3405467'`n" +
            "MATCH (s)-[:study_of_program]->(p)`n" +
            "MATCH (ss)<-[:sample_of_study_subject]-(samp)`n" +
            "MATCH (samp)<-[:file_of_sample]-(f)`n" +
            "MATCH (lp)<-[:file_of_laboratory_procedure]-(f)`n" +
            "RETURN COUNT(DISTINCT p) AS Programs,`n" +
            "COUNT(DISTINCT s) AS Arms,`n" +
            "COUNT(DISTINCT ss) AS Cases,`n" +
            "COUNT(DISTINCT samp) AS Samples,`n" +
            "COUNT(DISTINCT lp) AS Assays,`n" +
            "COUNT(DISTINCT f) AS Files`n"

# Same StatQuery text is used for CasesTab, SamplesTab and FilesTab rows.
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# The FilesTab row's wrapped query text now needs more vertical room.
$ws.Rows.Item(4).RowHeight = 255

# Move the saved selection/active cell to D12.
$ws.Range("D12").Select()
